$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host $ws.Name
$ws2 = $wb.Worksheets.Item(2)
Write-Host $ws2.Name
Write-Host $ws.Range("B1").Value
Write-Host $ws2.Range("A1").Value
